$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in product rows 8-11 (SKU, 平台仓, 重量, 采购单价, 头程费用)
$ws.Range("A8").Value = "效果器单块"
$ws.Range("B8").Value = "蓝舸"
$ws.Range("C8").Value = 0.18
$ws.Range("D8").Value = 38
$ws.Range("E8").Value = 1

$ws.Range("A9").Value = "调音器"
$ws.Range("B9").Value = "蓝舸"
$ws.Range("C9").Value = 0.2
$ws.Range("D9").Value = 10
$ws.Range("E9").Value = 1

$ws.Range("A10").Value = "NRXD68865"
$ws.Range("B10").Value = "拉美"
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 22.03
$ws.Range("E10").Value = 9.14

$ws.Range("A11").Value = "水槽"
$ws.Range("B11").Value = "蓝舸"
$ws.Range("C11").Value = 14
$ws.Range("D11").Value = 305
$ws.Range("E11").Value = 138

# Update selected cell on the active sheet
$ws.Range("D15").Select()
